$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: remove the F186 "i" marker cell (text -> cleared) ---
$ws.Range("F186").ClearContents()

# --- Step 2: convert "o" marker text cells in column F to numeric 1 ---
$ws.Range("F223").Value2 = 1
$ws.Range("F224").Value2 = 1
$ws.Range("F231").Value2 = 1
$ws.Range("F242").Value2 = 1
$ws.Range("F243").Value2 = 1
$ws.Range("F250").Value2 = 1
$ws.Range("F256").Value2 = 1

# --- Step 3: append new log rows 381-470, copying number formats from row 380 ---
$ws.Range("A380:B380").Copy()
$ws.Range("A381:B470").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# row 381
$ws.Range("A381").Value2 = 44618
$ws.Range("B381").Value2 = 0.38194444444444442
$ws.Range("C381").Value2 = 1
# row 382
$ws.Range("A382").Value2 = 44618
$ws.Range("B382").Value2 = 0.39861111111111108
$ws.Range("E382").Value2 = 1
# row 383
$ws.Range("A383").Value2 = 44618
$ws.Range("B383").Value2 = 0.46180555555555558
$ws.Range("D383").Value2 = 1
$ws.Range("F383").Value2 = 1
# row 384
$ws.Range("A384").Value2 = 44618
$ws.Range("B384").Value2 = 0.47430555555555554
$ws.Range("D384").Value2 = 1
$ws.Range("F384").Value2 = 1
# row 385
$ws.Range("A385").Value2 = 44618
$ws.Range("B385").Value2 = 0.62847222222222221
$ws.Range("C385").Value2 = 1
$ws.Range("F385").Value2 = 1
# row 386
$ws.Range("A386").Value2 = 44618
$ws.Range("B386").Value2 = 0.77083333333333337
$ws.Range("C386").Value2 = 1
# row 387
$ws.Range("A387").Value2 = 44618
$ws.Range("B387").Value2 = 0.79166666666666663
$ws.Range("E387").Value2 = 1
# row 388
$ws.Range("A388").Value2 = 44619
$ws.Range("B388").Value2 = 0.38194444444444442
$ws.Range("C388").Value2 = 1
# row 389
$ws.Range("A389").Value2 = 44619
$ws.Range("B389").Value2 = 0.40277777777777773
$ws.Range("C389").Value2 = 1
# row 390
$ws.Range("A390").Value2 = 44619
$ws.Range("B390").Value2 = 0.42569444444444443
$ws.Range("D390").Value2 = 1
$ws.Range("F390").Value2 = 1
# row 391
$ws.Range("A391").Value2 = 44619
$ws.Range("B391").Value2 = 0.45416666666666666
$ws.Range("C391").Value2 = 1
# row 392
$ws.Range("A392").Value2 = 44619
$ws.Range("B392").Value2 = 0.47500000000000003
$ws.Range("E392").Value2 = 1
# row 393
$ws.Range("A393").Value2 = 44619
$ws.Range("B393").Value2 = 0.87847222222222221
$ws.Range("C393").Value2 = 1
# row 394
$ws.Range("A394").Value2 = 44619
$ws.Range("B394").Value2 = 0.93402777777777779
$ws.Range("C394").Value2 = 1
# row 395
$ws.Range("A395").Value2 = 44619
$ws.Range("B395").Value2 = 0.95486111111111116
$ws.Range("D395").Value2 = 1
$ws.Range("F395").Value2 = 1
# row 396
$ws.Range("A396").Value2 = 44620
$ws.Range("B396").Value2 = 0.33333333333333331
$ws.Range("C396").Value2 = 1
# row 397
$ws.Range("A397").Value2 = 44620
$ws.Range("B397").Value2 = 0.35416666666666669
$ws.Range("E397").Value2 = 1
# row 398
$ws.Range("A398").Value2 = 44620
$ws.Range("B398").Value2 = 0.60069444444444442
$ws.Range("C398").Value2 = 1
# row 399
$ws.Range("A399").Value2 = 44620
$ws.Range("B399").Value2 = 0.81944444444444453
$ws.Range("E399").Value2 = 1
# row 400
$ws.Range("A400").Value2 = 44620
$ws.Range("B400").Value2 = 0.88888888888888884
$ws.Range("C400").Value2 = 0
# row 401
$ws.Range("A401").Value2 = 44620
$ws.Range("B401").Value2 = 0.91319444444444453
$ws.Range("C401").Value2 = 1
# row 402
$ws.Range("A402").Value2 = 44620
$ws.Range("B402").Value2 = 0.92361111111111116
$ws.Range("D402").Value2 = 1
# row 403
$ws.Range("A403").Value2 = 44621
$ws.Range("B403").Value2 = 0.35416666666666669
$ws.Range("C403").Value2 = 1
# row 404
$ws.Range("A404").Value2 = 44621
$ws.Range("B404").Value2 = 0.39930555555555558
$ws.Range("C404").Value2 = 0
# row 405
$ws.Range("A405").Value2 = 44621
$ws.Range("B405").Value2 = 0.40208333333333335
$ws.Range("D405").Value2 = 1
# row 406
$ws.Range("A406").Value2 = 44621
$ws.Range("B406").Value2 = 0.4055555555555555
$ws.Range("E406").Value2 = 1
# row 407
$ws.Range("A407").Value2 = 44621
$ws.Range("B407").Value2 = 0.71527777777777779
$ws.Range("C407").Value2 = 1
# row 408
$ws.Range("A408").Value2 = 44621
$ws.Range("B408").Value2 = 0.80902777777777779
$ws.Range("C408").Value2 = 1
# row 409
$ws.Range("A409").Value2 = 44621
$ws.Range("B409").Value2 = 0.91875000000000007
$ws.Range("D409").Value2 = 1
# row 410
$ws.Range("A410").Value2 = 44622
$ws.Range("B410").Value2 = 0.3611111111111111
$ws.Range("C410").Value2 = 1
# row 411
$ws.Range("A411").Value2 = 44622
$ws.Range("B411").Value2 = 0.37083333333333335
$ws.Range("E411").Value2 = 1
# row 412
$ws.Range("A412").Value2 = 44622
$ws.Range("B412").Value2 = 0.4826388888888889
$ws.Range("C412").Value2 = 1
$ws.Range("F412").Value2 = 1
# row 413
$ws.Range("A413").Value2 = 44622
$ws.Range("B413").Value2 = 0.48749999999999999
$ws.Range("D413").Value2 = 1
$ws.Range("F413").Value2 = 1
# row 414
$ws.Range("A414").Value2 = 44622
$ws.Range("B414").Value2 = 0.63541666666666663
$ws.Range("C414").Value2 = 1
# row 415
$ws.Range("A415").Value2 = 44622
$ws.Range("B415").Value2 = 0.77083333333333337
$ws.Range("E415").Value2 = 1
# row 416
$ws.Range("A416").Value2 = 44622
$ws.Range("B416").Value2 = 0.91875000000000007
$ws.Range("C416").Value2 = 1
# row 417
$ws.Range("A417").Value2 = 44622
$ws.Range("B417").Value2 = 0.9277777777777777
$ws.Range("D417").Value2 = 1
# row 418
$ws.Range("A418").Value2 = 44622
$ws.Range("B418").Value2 = 0.35138888888888892
$ws.Range("C418").Value2 = 1
# row 419
$ws.Range("A419").Value2 = 44622
$ws.Range("B419").Value2 = 0.35486111111111113
$ws.Range("E419").Value2 = 1
# row 420
$ws.Range("A420").Value2 = 44622
$ws.Range("B420").Value2 = 0.4236111111111111
$ws.Range("C420").Value2 = 1
# row 421
$ws.Range("A421").Value2 = 44622
$ws.Range("B421").Value2 = 0.42777777777777781
$ws.Range("D421").Value2 = 1
# row 422
$ws.Range("A422").Value2 = 44622
$ws.Range("B422").Value2 = 0.55208333333333337
$ws.Range("C422").Value2 = 1
# row 423
$ws.Range("A423").Value2 = 44622
$ws.Range("B423").Value2 = 0.76736111111111116
$ws.Range("C423").Value2 = 1
$ws.Range("F423").Value2 = 1
# row 424
$ws.Range("A424").Value2 = 44623
$ws.Range("B424").Value2 = 0.34027777777777773
$ws.Range("C424").Value2 = 1
# row 425
$ws.Range("A425").Value2 = 44623
$ws.Range("B425").Value2 = 0.37916666666666665
$ws.Range("D425").Value2 = 1
$ws.Range("F425").Value2 = 1
# row 426
$ws.Range("A426").Value2 = 44623
$ws.Range("B426").Value2 = 0.3888888888888889
$ws.Range("E426").Value2 = 1
# row 427
$ws.Range("A427").Value2 = 44623
$ws.Range("B427").Value2 = 0.70833333333333337
$ws.Range("C427").Value2 = 1
# row 428
$ws.Range("A428").Value2 = 44623
$ws.Range("B428").Value2 = 0.79236111111111107
$ws.Range("C428").Value2 = 1
# row 429
$ws.Range("A429").Value2 = 44623
$ws.Range("B429").Value2 = 0.80555555555555547
$ws.Range("E429").Value2 = 1
# row 430
$ws.Range("A430").Value2 = 44623
$ws.Range("B430").Value2 = 0.90555555555555556
$ws.Range("C430").Value2 = 1
# row 431
$ws.Range("A431").Value2 = 44623
$ws.Range("B431").Value2 = 0.91180555555555554
$ws.Range("D431").Value2 = 1
# row 432
$ws.Range("A432").Value2 = 44624
$ws.Range("B432").Value2 = 0.39583333333333331
$ws.Range("C432").Value2 = 1
# row 433
$ws.Range("A433").Value2 = 44624
$ws.Range("B433").Value2 = 0.40625
$ws.Range("E433").Value2 = 1
# row 434
$ws.Range("A434").Value2 = 44624
$ws.Range("B434").Value2 = 0.47291666666666665
$ws.Range("C434").Value2 = 0
# row 435
$ws.Range("A435").Value2 = 44624
$ws.Range("B435").Value2 = 0.47638888888888892
$ws.Range("D435").Value2 = 1
# row 436
$ws.Range("A436").Value2 = 44624
$ws.Range("B436").Value2 = 0.59583333333333333
$ws.Range("C436").Value2 = 1
$ws.Range("F436").Value2 = 1
# row 437
$ws.Range("A437").Value2 = 44624
$ws.Range("B437").Value2 = 0.60069444444444442
$ws.Range("D437").Value2 = 1
$ws.Range("F437").Value2 = 1
# row 438
$ws.Range("A438").Value2 = 44624
$ws.Range("B438").Value2 = 0.68194444444444446
$ws.Range("C438").Value2 = 1
# row 439
$ws.Range("A439").Value2 = 44624
$ws.Range("B439").Value2 = 0.80486111111111114
$ws.Range("E439").Value2 = 1
# row 440
$ws.Range("A440").Value2 = 44624
$ws.Range("B440").Value2 = 0.90486111111111101
$ws.Range("C440").Value2 = 1
# row 441
$ws.Range("A441").Value2 = 44625
$ws.Range("B441").Value2 = 0.39583333333333331
$ws.Range("C441").Value2 = 1
# row 442
$ws.Range("A442").Value2 = 44625
$ws.Range("B442").Value2 = 0.40625
$ws.Range("E442").Value2 = 1
# row 443
$ws.Range("A443").Value2 = 44625
$ws.Range("B443").Value2 = 0.47291666666666665
$ws.Range("C443").Value2 = 0
# row 444
$ws.Range("A444").Value2 = 44625
$ws.Range("B444").Value2 = 0.47638888888888892
$ws.Range("D444").Value2 = 1
# row 445
$ws.Range("A445").Value2 = 44625
$ws.Range("B445").Value2 = 0.59583333333333333
$ws.Range("C445").Value2 = 1
$ws.Range("F445").Value2 = 1
# row 446
$ws.Range("A446").Value2 = 44625
$ws.Range("B446").Value2 = 0.60069444444444442
$ws.Range("D446").Value2 = 1
$ws.Range("F446").Value2 = 1
# row 447
$ws.Range("A447").Value2 = 44625
$ws.Range("B447").Value2 = 0.68194444444444446
$ws.Range("C447").Value2 = 1
$ws.Range("F447").Value2 = 1
# row 448
$ws.Range("A448").Value2 = 44625
$ws.Range("B448").Value2 = 0.80486111111111114
$ws.Range("E448").Value2 = 1
# row 449
$ws.Range("A449").Value2 = 44625
$ws.Range("B449").Value2 = 0.90486111111111101
$ws.Range("C449").Value2 = 1
# row 450
$ws.Range("A450").Value2 = 44626
$ws.Range("B450").Value2 = 0.3923611111111111
$ws.Range("C450").Value2 = 1
# row 451
$ws.Range("A451").Value2 = 44626
$ws.Range("B451").Value2 = 0.40625
$ws.Range("E451").Value2 = 1
# row 452
$ws.Range("A452").Value2 = 44626
$ws.Range("B452").Value2 = 0.51944444444444449
$ws.Range("C452").Value2 = 1
$ws.Range("F452").Value2 = 1
# row 453
$ws.Range("A453").Value2 = 44626
$ws.Range("B453").Value2 = 0.52222222222222225
$ws.Range("D453").Value2 = 1
$ws.Range("F453").Value2 = 1
# row 454
$ws.Range("A454").Value2 = 44626
$ws.Range("B454").Value2 = 0.52638888888888891
$ws.Range("D454").Value2 = 1
$ws.Range("F454").Value2 = 1
# row 455
$ws.Range("A455").Value2 = 44626
$ws.Range("B455").Value2 = 0.52986111111111112
$ws.Range("D455").Value2 = 1
$ws.Range("F455").Value2 = 1
# row 456
$ws.Range("A456").Value2 = 44626
$ws.Range("B456").Value2 = 0.5805555555555556
$ws.Range("C456").Value2 = 1
$ws.Range("F456").Value2 = 1
# row 457
$ws.Range("A457").Value2 = 44626
$ws.Range("B457").Value2 = 0.58194444444444449
$ws.Range("C457").Value2 = 1
$ws.Range("F457").Value2 = 1
# row 458
$ws.Range("A458").Value2 = 44626
$ws.Range("B458").Value2 = 0.74305555555555547
$ws.Range("C458").Value2 = 1
# row 459
$ws.Range("A459").Value2 = 44626
$ws.Range("B459").Value2 = 0.79513888888888884
$ws.Range("C459").Value2 = 1
$ws.Range("F459").Value2 = 1
# row 460
$ws.Range("A460").Value2 = 44626
$ws.Range("B460").Value2 = 0.80902777777777779
$ws.Range("E460").Value2 = 1
# row 461
$ws.Range("A461").Value2 = 44626
$ws.Range("B461").Value2 = 0.91875000000000007
$ws.Range("C461").Value2 = 1
# row 462
$ws.Range("A462").Value2 = 44627
$ws.Range("B462").Value2 = 0.31944444444444448
$ws.Range("C462").Value2 = 1
# row 463
$ws.Range("A463").Value2 = 44627
$ws.Range("B463").Value2 = 0.33611111111111108
$ws.Range("C463").Value2 = 1
$ws.Range("F463").Value2 = 1
# row 464
$ws.Range("A464").Value2 = 44627
$ws.Range("B464").Value2 = 0.3611111111111111
$ws.Range("E464").Value2 = 1
# row 465
$ws.Range("A465").Value2 = 44627
$ws.Range("B465").Value2 = 0.45416666666666666
$ws.Range("C465").Value2 = 1
# row 466
$ws.Range("A466").Value2 = 44627
$ws.Range("B466").Value2 = 0.60625000000000007
$ws.Range("C466").Value2 = 1
# row 467
$ws.Range("A467").Value2 = 44627
$ws.Range("B467").Value2 = 0.75347222222222221
$ws.Range("C467").Value2 = 1
# row 468
$ws.Range("A468").Value2 = 44627
$ws.Range("B468").Value2 = 0.77083333333333337
$ws.Range("C468").Value2 = 1
$ws.Range("F468").Value2 = 1
# row 469
$ws.Range("A469").Value2 = 44627
$ws.Range("B469").Value2 = 0.77777777777777779
$ws.Range("D469").Value2 = 1
$ws.Range("F469").Value2 = 1
# row 470
$ws.Range("A470").Value2 = 44627
$ws.Range("B470").Value2 = 0.78472222222222221
$ws.Range("D470").Value2 = 1
$ws.Range("F470").Value2 = 1

# --- Step 4: update view selection to A2 (matches frozen-pane view after edit) ---
$ws.Range("A2").Select()
